$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Datos actualizados a 18 de Abril de 2020 a las 20:52
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 20:52"

# Row 7: Francia
$ws.Range("A7").Value = "Francia"
$ws.Range("B7").Value = 151793
$ws.Range("C7").Value = 3824
$ws.Range("D7").Value = 35983
$ws.Range("E7").Value = 96487
$ws.Range("F7").Value = 5833
$ws.Range("G7").Value = 642
$ws.Range("H7").Value = 19323

# Row 16: Canada
$ws.Range("A16").Value = "Canada"
$ws.Range("B16").Value = 33180
$ws.Range("C16").Value = 1253
$ws.Range("D16").Value = 11141
$ws.Range("E16").Value = 20573
$ws.Range("F16").Value = 557
$ws.Range("G16").Value = 156
$ws.Range("H16").Value = 1466

# Row 18: Suiza
$ws.Range("A18").Value = "Suiza"
$ws.Range("B18").Value = 27404
$ws.Range("C18").Value = 326
$ws.Range("D18").Value = 17100
$ws.Range("E18").Value = 8938
$ws.Range("F18").Value = 386
$ws.Range("G18").Value = 39
$ws.Range("H18").Value = 1366

# Row 23: Peru
$ws.Range("A23").Value = "Peru"
$ws.Range("B23").Value = 14420
$ws.Range("C23").Value = 931
$ws.Range("D23").Value = 6684
$ws.Range("E23").Value = 7388
$ws.Range("F23").Value = 117
$ws.Range("G23").Value = 48
$ws.Range("H23").Value = 348

# Row 64: Barein
$ws.Range("A64").Value = "Barein"
$ws.Range("B64").Value = 1773
$ws.Range("C64").Value = 33
$ws.Range("D64").Value = 741
$ws.Range("E64").Value = 1025
$ws.Range("F64").Value = 3
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 7

# Row 96: Burkina Faso
$ws.Range("A96").Value = "Burkina Faso"
$ws.Range("B96").Value = 565
$ws.Range("C96").Value = 8
$ws.Range("D96").Value = 321
$ws.Range("E96").Value = 208
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 36

# Row 118: Sri Lanka
$ws.Range("A118").Value = "Sri Lanka"
$ws.Range("B118").Value = 254
$ws.Range("C118").Value = 10
$ws.Range("D118").Value = 86
$ws.Range("E118").Value = 161
$ws.Range("F118").Value = 1
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 7

# Row 119: Mayotte
$ws.Range("A119").Value = "Mayotte"
$ws.Range("B119").Value = 254
$ws.Range("C119").Value = 9
$ws.Range("D119").Value = 117
$ws.Range("E119").Value = 133
$ws.Range("F119").Value = 6
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 4

# Row 163: Eritrea
$ws.Range("A163").Value = "Eritrea"
$ws.Range("B163").Value = 39
$ws.Range("C163").Value = 4
$ws.Range("D163").Value = 0
$ws.Range("E163").Value = 39
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 0

# Row 164: Puerto Rico
$ws.Range("A164").Value = "Puerto Rico"
$ws.Range("B164").Value = 39
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 1
$ws.Range("E164").Value = 36
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 2

# Row 165: Siria
$ws.Range("A165").Value = "Siria"
$ws.Range("B165").Value = 38
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 5
$ws.Range("E165").Value = 31
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 2

# Row 166: San Martin (Parte Francesa)
$ws.Range("A166").Value = "San Martin (Parte Francesa)"
$ws.Range("B166").Value = 37
$ws.Range("C166").Value = 2
$ws.Range("D166").Value = 19
$ws.Range("E166").Value = 16
$ws.Range("F166").Value = 5
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 2

# Row 176: Angola
$ws.Range("A176").Value = "Angola"
$ws.Range("B176").Value = 24
$ws.Range("C176").Value = 5
$ws.Range("D176").Value = 6
$ws.Range("E176").Value = 16
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 2

# Row 177: Antigua y Barbuda
$ws.Range("A177").Value = "Antigua y Barbuda"
$ws.Range("B177").Value = 23
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 3
$ws.Range("E177").Value = 17
$ws.Range("F177").Value = 1
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 3

# Row 178: Suazilandia
$ws.Range("A178").Value = "Suazilandia"
$ws.Range("B178").Value = 22
$ws.Range("C178").Value = 6
$ws.Range("D178").Value = 8
$ws.Range("E178").Value = 13
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 1

# Row 179: Laos
$ws.Range("A179").Value = "Laos"
$ws.Range("B179").Value = 19
$ws.Range("C179").Value = 0
$ws.Range("D179").Value = 2
$ws.Range("E179").Value = 17
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 0

# Row 215: San Pedro y Miquelon
$ws.Range("A215").Value = "San Pedro y Miquelon"
$ws.Range("B215").Value = 1
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 0
$ws.Range("E215").Value = 1
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

# Row 216: Yemen
$ws.Range("A216").Value = "Yemen"
$ws.Range("B216").Value = 1
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 0
$ws.Range("E216").Value = 1
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0

